# Update "想去人数" (number of people wanting to go) counts that changed
# between scrapes, on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitionUpdates = @{
    "F2"  = 824
    "F3"  = 24
    "F9"  = 538
    "F12" = 13312
    "F13" = 165
    "F16" = 5477
    "F17" = 5569
    "F18" = 40
}

foreach ($addr in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range($addr).Value = $exhibitionUpdates[$addr]
}

# Sheet "全部类型": row -> new F value
$allUpdates = @{
    "F2"  = 824
    "F10" = 24
    "F31" = 538
    "F34" = 13312
    "F35" = 165
    "F39" = 5477
    "F40" = 5569
    "F41" = 40
}

foreach ($addr in $allUpdates.Keys) {
    $sheetAll.Range($addr).Value = $allUpdates[$addr]
}

$wb.Save()
